$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2026-02-18 Wednesday" "2026-02-19 Thursday"

Replace-Text "90×21=1890" "41×60=2460"
Replace-Text "27×14=378" "41×89=3649"
Replace-Text "21×99=2079" "77×11=847"
Replace-Text "63×37=2331" "48×81=3888"
Replace-Text "57×88=5016" "21×81=1701"

Replace-Text "76×47=3572" "68×57=3876"
Replace-Text "47×11=517" "94×32=3008"
Replace-Text "69×95=6555" "77×69=5313"
Replace-Text "40×50=2000" "21×58=1218"
Replace-Text "95×33=3135" "74×22=1628"

Replace-Text "12×64=768" "45×95=4275"
Replace-Text "80×65=5200" "47×78=3666"
Replace-Text "50×15=750" "11×74=814"
Replace-Text "64×65=4160" "89×18=1602"
Replace-Text "74×59=4366" "61×67=4087"

Replace-Text "20×47=940" "73×19=1387"
Replace-Text "65×46=2990" "67×88=5896"
Replace-Text "50×31=1550" "17×71=1207"
Replace-Text "60×73=4380" "36×96=3456"
Replace-Text "98×81=7938" "28×39=1092"

Replace-Text "57×69=3933" "78×20=1560"
Replace-Text "69×35=2415" "53×66=3498"
Replace-Text "20×24=480" "33×95=3135"
Replace-Text "30×67=2010" "41×77=3157"
Replace-Text "46×94=4324" "62×58=3596"
